$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header capitalization: "Car price" -> "Car Price"
$ws.Range("B1").Value = "Car Price"

# The Name/Price columns had been swapped (and the data shuffled) in the
# original file -- restore the correct data with Car Name in column A and
# Car Price in column B for every data row.
$ws.Range("A2").Value = "Honda Brio 2011-2013 Honda Brio S MT"
$ws.Range("B2").Value = "Rs. 3.16 Lakh"

$ws.Range("A3").Value = "Honda Amaze 2013-2016 Honda Amaze S AT i-Vtech"
$ws.Range("B3").Value = "Rs. 3.90 Lakh"

$ws.Range("A4").Value = "Honda Brio 2011-2013 Honda Brio S MT"
$ws.Range("B4").Value = "Rs. 3.05 Lakh"

$ws.Range("A5").Value = "Honda Brio 2013-2016 Honda Brio S MT"
$ws.Range("B5").Value = "Rs. 3.50 Lakh"

$ws.Range("A6").Value = "Honda Brio 2013-2016 Honda Brio S MT"
$ws.Range("B6").Value = "Rs. 3.60 Lakh"

$ws.Range("A7").Value = "Honda City 2008-2011 Honda City 1.5 S MT"
$ws.Range("B7").Value = "Rs. 4.00 Lakh"

$ws.Range("A8").Value = "Honda Jazz 2014-2020 Honda Jazz 1.2 SV i VTEC"
$ws.Range("B8").Value = "Rs. 3.85 Lakh"

$ws.Range("A9").Value = "Honda City ZX 2005-2008 Honda City GXi"
$ws.Range("B9").Value = "Rs. 2.10 Lakh"

$ws.Range("A10").Value = "Honda Amaze 2013-2016 Honda Amaze VX i-Vtech"
$ws.Range("B10").Value = "Rs. 4.00 Lakh"

$ws.Range("A11").Value = "Honda Civic 2010-2013 Honda Civic 1.8 V MT"
$ws.Range("B11").Value = "Rs. 3.00 Lakh"

$ws.Range("A12").Value = "Honda Civic 2010-2013 Honda Civic 1.8 V MT"
$ws.Range("B12").Value = "Rs. 3.00 Lakh"

$ws.Range("A13").Value = "Honda Brio 2011-2013 Honda Brio E MT"
$ws.Range("B13").Value = "Rs. 1.90 Lakh"
